$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 93.8
$ws.Range("I4").Value = 93.8
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 93.8
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 20.2

$ws.Range("H15").Value = 124.15
$ws.Range("I15").Value = 124.15
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 372.45
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -203.45

$ws.Range("H28").Value = 570.61536
$ws.Range("I28").Value = 598.86365
$ws.Range("J28").Value = 415.25
$ws.Range("K28").Value = 598.86365
$ws.Range("L28").Value = 415.25
$ws.Range("M28").Value = -113.86365
$ws.Range("N28").Value = -1385.25

$ws.Range("H32").Value = 13940036
$ws.Range("I32").Value = 780
$ws.Range("J32").Value = 23232872
$ws.Range("K32").Value = 780
$ws.Range("L32").Value = 23232872
$ws.Range("M32").Value = -454
$ws.Range("N32").Value = -23233524

$ws.Range("H40").Value = 1116.7778
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1123.6471
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1123.6471
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1473.6471

$ws.Range("H62").Value = 10805.4375
$ws.Range("I62").Value = 14786.777
$ws.Range("J62").Value = 5686.5713
$ws.Range("K62").Value = 14786.777
$ws.Range("L62").Value = 5686.5713
$ws.Range("M62").Value = -14162.777
$ws.Range("N62").Value = -6934.5713

$ws.Range("H65").Value = 10805.4375
$ws.Range("I65").Value = 14786.777
$ws.Range("J65").Value = 5686.5713
$ws.Range("K65").Value = 73933.88499999999
$ws.Range("L65").Value = 28432.8565
$ws.Range("M65").Value = -70813.88499999999
$ws.Range("N65").Value = -34672.85649999999

$ws.Range("H107").Value = 702.6842
$ws.Range("I107").Value = 774
$ws.Range("J107").Value = 435.25
$ws.Range("K107").Value = 774
$ws.Range("L107").Value = 435.25
$ws.Range("M107").Value = 1146
$ws.Range("N107").Value = -4275.25

$ws.Range("H116").Value = 165337.31
$ws.Range("I116").Value = 194080.45
$ws.Range("J116").Value = 7250
$ws.Range("K116").Value = 194080.45
$ws.Range("L116").Value = 7250
$ws.Range("M116").Value = -190638.45
$ws.Range("N116").Value = -14134

$ws.Range("H125").Value = 1045.0769
$ws.Range("I125").Value = 1446.6666
$ws.Range("J125").Value = 700.8570999999999
$ws.Range("K125").Value = 13019.9994
$ws.Range("L125").Value = 6307.7139
$ws.Range("M125").Value = -10559.9994
$ws.Range("N125").Value = -11227.7139

$ws.Range("H131").Value = 2294.6875
$ws.Range("I131").Value = 1492.9166
$ws.Range("J131").Value = 4700
$ws.Range("K131").Value = 4478.7498
$ws.Range("L131").Value = 14100
$ws.Range("M131").Value = 561.2502000000004
$ws.Range("N131").Value = -24180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3006.6897
$ws.Range("I61").Value = 1966.8572
$ws.Range("J61").Value = 3977.2
$ws.Range("K61").Value = 1966.8572
$ws.Range("L61").Value = 3977.2
$ws.Range("M61").Value = -1754.8572
$ws.Range("N61").Value = -4401.2

$ws.Range("H97").Value = 739.56525
$ws.Range("I97").Value = 702.9231
$ws.Range("J97").Value = 943.7143
$ws.Range("K97").Value = 702.9231
$ws.Range("L97").Value = 943.7143
$ws.Range("M97").Value = -206.9231
$ws.Range("N97").Value = -1935.7143

$ws.Range("H110").Value = 1214.3939
$ws.Range("I110").Value = 1138.4166
$ws.Range("J110").Value = 1417
$ws.Range("K110").Value = 1138.4166
$ws.Range("L110").Value = 1417
$ws.Range("M110").Value = 906.5834
$ws.Range("N110").Value = -5507

$ws.Range("H132").Value = 3788.0408
$ws.Range("I132").Value = 3884.963
$ws.Range("J132").Value = 3669.0908
$ws.Range("K132").Value = 11654.889
$ws.Range("L132").Value = 11007.2724
$ws.Range("M132").Value = -9124.889000000001
$ws.Range("N132").Value = -16067.2724

$ws.Range("H136").Value = 3006.6897
$ws.Range("I136").Value = 1966.8572
$ws.Range("J136").Value = 3977.2
$ws.Range("K136").Value = 5900.571599999999
$ws.Range("L136").Value = 11931.6
$ws.Range("M136").Value = -3350.571599999999
$ws.Range("N136").Value = -17031.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 498
$ws.Range("I22").Value = 447.5
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 447.5
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -274.5
$ws.Range("N22").Value = -1046

$ws.Range("H80").Value = 315.64
$ws.Range("I80").Value = 177.375
$ws.Range("J80").Value = 380.70587
$ws.Range("K80").Value = 177.375
$ws.Range("L80").Value = 380.70587
$ws.Range("M80").Value = 820.625
$ws.Range("N80").Value = -2376.70587

$ws.Range("H83").Value = 315.64
$ws.Range("I83").Value = 177.375
$ws.Range("J83").Value = 380.70587
$ws.Range("K83").Value = 886.875
$ws.Range("L83").Value = 1903.52935
$ws.Range("M83").Value = 4105.125
$ws.Range("N83").Value = -11887.52935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3259.1
$ws.Range("I31").Value = 2380.8572
$ws.Range("J31").Value = 4027.5625
$ws.Range("K31").Value = 2380.8572
$ws.Range("L31").Value = 4027.5625
$ws.Range("M31").Value = -2085.8572
$ws.Range("N31").Value = -4617.5625

$ws.Range("H34").Value = 3259.1
$ws.Range("I34").Value = 2380.8572
$ws.Range("J34").Value = 4027.5625
$ws.Range("K34").Value = 2380.8572
$ws.Range("L34").Value = 4027.5625
$ws.Range("M34").Value = -2178.8572
$ws.Range("N34").Value = -4431.5625

$ws.Range("H105").Value = 621.8
$ws.Range("I105").Value = 609.0909
$ws.Range("J105").Value = 715
$ws.Range("K105").Value = 609.0909
$ws.Range("L105").Value = 715
$ws.Range("M105").Value = 1137.9091
$ws.Range("N105").Value = -4209

$ws.Range("H122").Value = 918.46344
$ws.Range("I122").Value = 754.0454999999999
$ws.Range("J122").Value = 1108.8422
$ws.Range("K122").Value = 2262.1365
$ws.Range("L122").Value = 3326.5266
$ws.Range("M122").Value = 187.8635000000004
$ws.Range("N122").Value = -8226.526600000001

$ws.Range("H132").Value = 3036.1072
$ws.Range("I132").Value = 1299.9166
$ws.Range("J132").Value = 4338.25
$ws.Range("K132").Value = 3899.7498
$ws.Range("L132").Value = 13014.75
$ws.Range("M132").Value = -1369.7498
$ws.Range("N132").Value = -18074.75

$ws.Range("H134").Value = 2125.4333
$ws.Range("I134").Value = 1267.0769
$ws.Range("J134").Value = 2781.8235
$ws.Range("K134").Value = 3801.2307
$ws.Range("L134").Value = 8345.470499999999
$ws.Range("M134").Value = -1266.2307
$ws.Range("N134").Value = -13415.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 40.22222
$ws.Range("I38").Value = 37.384617
$ws.Range("J38").Value = 47.6
$ws.Range("K38").Value = 112.153851
$ws.Range("L38").Value = 142.8
$ws.Range("M38").Value = 234.846149
$ws.Range("N38").Value = -836.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4948.524
$ws.Range("I102").Value = 5206.5
$ws.Range("J102").Value = 3400.6667
$ws.Range("K102").Value = 5206.5
$ws.Range("L102").Value = 3400.6667
$ws.Range("M102").Value = -3584.5
$ws.Range("N102").Value = -6644.6667

$ws.Range("H122").Value = 1616.6666
$ws.Range("I122").Value = 1466.6666
$ws.Range("J122").Value = 1766.6666
$ws.Range("K122").Value = 4399.9998
$ws.Range("L122").Value = 5299.9998
$ws.Range("M122").Value = -1949.9998
$ws.Range("N122").Value = -10199.9998

$ws.Range("H132").Value = 2915.5
$ws.Range("I132").Value = 2856.4167
$ws.Range("J132").Value = 3021.85
$ws.Range("K132").Value = 8569.250100000001
$ws.Range("L132").Value = 9065.549999999999
$ws.Range("M132").Value = -6039.250100000001
$ws.Range("N132").Value = -14125.55

$ws.Range("H136").Value = 8489.6
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8489.6
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 25468.8
$ws.Range("N136").Value = -30568.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3950.5
$ws.Range("I40").Value = 3940.8
$ws.Range("J40").Value = 3966.6667
$ws.Range("K40").Value = 3940.8
$ws.Range("L40").Value = 3966.6667
$ws.Range("M40").Value = -3804.8
$ws.Range("N40").Value = -4238.6667

$ws.Range("H122").Value = 3866.6667
$ws.Range("I122").Value = 3100
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 9300
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -6850
$ws.Range("N122").Value = -17650

$ws.Range("H136").Value = 6376.7417
$ws.Range("I136").Value = 4533.25
$ws.Range("J136").Value = 7541.0527
$ws.Range("K136").Value = 13599.75
$ws.Range("L136").Value = 22623.1581
$ws.Range("M136").Value = -11049.75
$ws.Range("N136").Value = -27723.1581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 275.33334
$ws.Range("I113").Value = 220.85715
$ws.Range("J113").Value = 580.4
$ws.Range("K113").Value = 662.5714499999999
$ws.Range("L113").Value = 1741.2
$ws.Range("M113").Value = 1507.42855
$ws.Range("N113").Value = -6081.2

$ws.Range("H126").Value = 3290.25
$ws.Range("I126").Value = 4752.6665
$ws.Range("J126").Value = 2412.8
$ws.Range("K126").Value = 14257.9995
$ws.Range("L126").Value = 7238.400000000001
$ws.Range("M126").Value = -11787.9995
$ws.Range("N126").Value = -12178.4

$ws.Range("H132").Value = 20258.89
$ws.Range("I132").Value = 32520.844
$ws.Range("J132").Value = 3198.7827
$ws.Range("K132").Value = 97562.53200000001
$ws.Range("L132").Value = 9596.348100000001
$ws.Range("M132").Value = -95032.53200000001
$ws.Range("N132").Value = -14656.3481

$ws.Range("H136").Value = 25643296
$ws.Range("I136").Value = 52633444
$ws.Range("J136").Value = 2655
$ws.Range("K136").Value = 157900332
$ws.Range("L136").Value = 7965
$ws.Range("M136").Value = -157897782
$ws.Range("N136").Value = -13065
